# Fruta / hortaliza, semanal
# The weekly refresh re-pulls the source rows, which results in the
# existing data rows (2-7) being re-ordered. Column A,B,C,E,F,G,H,I,J,K,Q,T
# stay identical for every row, so only columns D,L,M,N,O,P,R,S actually move
# between rows. Capture the "before" values for those columns, then write
# them back out under the new row mapping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by letter) whose values get shuffled between rows 2-7.
$cols = @("D", "L", "M", "N", "O", "P", "R", "S")

# Snapshot the current values for rows 2-7 before we start overwriting.
# NOTE: use Value2 (not Value) - Value round-trips incorrectly in this
# runtime when read back into a variable and re-assigned.
$snapshot = @{}
for ($r = 2; $r -le 7; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# New row r gets the old values that used to live in row $mapping[r].
$mapping = @{
    2 = 6
    3 = 7
    4 = 2
    5 = 4
    6 = 3
    7 = 5
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $srcVals[$col]
    }
}
